$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: moonshotai/kimi-k2-instruct-0905
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 0.001
$ws.Range("K10").Value = 469
$ws.Range("L10").Value = 0.001563333333333333

# Row 11: openai/gpt-oss-120b
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0.001
$ws.Range("K11").Value = 469
$ws.Range("L11").Value = 0.002345
